# Weekly price-sheet update: prepend a new week's record.
#
# The sheet is a flat, date-descending log of weekly price observations
# (one sheet, header in row 1, data in rows 2..343). The edit inserts a
# brand-new observation as the new row 244 (pushing the former rows
# 244..343 down to 245..344) — i.e. a new week's data point was recorded,
# duplicating the previous top-of-block record's metrics but dated 8
# weeks later (serial date 45027 instead of 44971).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shift rows 244:343 down to 245:344.
$ws.Rows.Item(244).Insert()

# Populate the newly-opened row 244 with the new observation.
$ws.Range("A244").Value = 4
$ws.Range("B244").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C244").Value = 'Los Lagos'
$ws.Range("D244").Value = 45027
$ws.Range("E244").Value = 10
$ws.Range("F244").Value = 100112039
$ws.Range("G244").Value = 'Ciboulette'
$ws.Range("H244").Value = 'Sin especificar'
$ws.Range("I244").Value = 'Primera'
$ws.Range("J244").Value = 240
$ws.Range("K244").Value = 3500
$ws.Range("L244").Value = 3500
$ws.Range("M244").Value = 3500
$ws.Range("N244").Value = '$/docena de atados'
$ws.Range("O244").Value = 'Región Metropolitana'
$ws.Range("P244").Value = 1167
$ws.Range("Q244").Value = 3
$ws.Range("R244").Value = 'Hortaliza'
